$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header "Sector" -> "Section"
$ws.Range("A1").Value = "Section"

# Update Sector/Section value (A2)
$ws.Range("A2").Value = "Полимерные компаунды"

# Update Description value (B2)
$ws.Range("B2").Value = "Данный патент описывает полимерную композицию на основе полиолефина, которая может использоваться для производства мембран и других изделий. Композиция содержит термопластичный полимер, произведенный в реакторе (rTPO) и по меньшей мере один линейный полиэтилен низкой плотности (LLDPE), что позволяет улучшить стабильность размеров и силу термического расширения по сравнению с композицией, содержащей только rTPO, без влияния на другие свойства. Композиция также может содержать пакет добавок и/или пакет наполнителей. "
